$wb = $excel.ActiveWorkbook

$ws1 = $wb.Worksheets.Item("ALC")
$ws2 = $wb.Worksheets.Item("ARM")
$ws3 = $wb.Worksheets.Item("BSM")
$ws4 = $wb.Worksheets.Item("CRP")
$ws5 = $wb.Worksheets.Item("CUL")
$ws6 = $wb.Worksheets.Item("GSM")
$ws7 = $wb.Worksheets.Item("LTW")
$ws8 = $wb.Worksheets.Item("WVR")

# ALC row 17
$ws1.Range("H17").Value = 1599.1842
$ws1.Range("I17").Value = 1150
$ws1.Range("J17").Value = 1624.1389
$ws1.Range("K17").Value = 3450
$ws1.Range("L17").Value = 4872.4167
$ws1.Range("M17").Value = -3282
$ws1.Range("N17").Value = -5208.4167

# ALC row 51
$ws1.Range("H51").Value = 6025.3438
$ws1.Range("I51").Value = 3745.818
$ws1.Range("J51").Value = 7219.381
$ws1.Range("K51").Value = 3745.818
$ws1.Range("L51").Value = 7219.381
$ws1.Range("M51").Value = -3261.818
$ws1.Range("N51").Value = -8187.381

# ALC row 55
$ws1.Range("H55").Value = 227.4
$ws1.Range("J55").Value = 289.6
$ws1.Range("L55").Value = 289.6
$ws1.Range("N55").Value = -717.6

# ALC row 98
$ws1.Range("H98").Value = 1736.8572
$ws1.Range("I98").Value = 1717.4375
$ws1.Range("K98").Value = 1717.4375
$ws1.Range("M98").Value = -219.4375

# ALC row 107
$ws1.Range("H107").Value = 20834306
$ws1.Range("I107").Value = 22223246
$ws1.Range("K107").Value = 22223246
$ws1.Range("M107").Value = -22221326

# ALC row 122
$ws1.Range("H122").Value = 1736.8572
$ws1.Range("I122").Value = 1717.4375
$ws1.Range("K122").Value = 5152.3125
$ws1.Range("M122").Value = -2702.3125

# ALC row 131
$ws1.Range("H131").Value = 3905.2
$ws1.Range("I131").Value = 3744.8462
$ws1.Range("K131").Value = 11234.5386
$ws1.Range("M131").Value = -6194.5386

# ALC row 137
$ws1.Range("H137").Value = 45447.953
$ws1.Range("I137").Value = 90193.89999999999
$ws1.Range("K137").Value = 270581.7
$ws1.Range("M137").Value = -268031.7

# ALC row 141
$ws1.Range("H141").Value = 2256.077
$ws1.Range("I141").Value = 1977.9
$ws1.Range("K141").Value = 5933.700000000001
$ws1.Range("M141").Value = -753.7000000000007

# ARM row 43
$ws2.Range("H43").Value = 34123.668
$ws2.Range("J43").Value = 34123.668
$ws2.Range("L43").Value = 34123.668
$ws2.Range("N43").Value = -34749.668

# ARM row 122
$ws2.Range("H122").Value = 1491368.4
$ws2.Range("I122").Value = 3373.111
$ws2.Range("J122").Value = 4169759.8
$ws2.Range("K122").Value = 10119.333
$ws2.Range("L122").Value = 12509279.4
$ws2.Range("M122").Value = -7669.332999999999
$ws2.Range("N122").Value = -12514179.4

# ARM row 132
$ws2.Range("H132").Value = 2390.48
$ws2.Range("I132").Value = 1492.5294
$ws2.Range("J132").Value = 4298.625
$ws2.Range("K132").Value = 4477.5882
$ws2.Range("L132").Value = 12895.875
$ws2.Range("M132").Value = -1947.5882
$ws2.Range("N132").Value = -17955.875

# BSM row 134
$ws3.Range("H134").Value = 5437.6
$ws3.Range("I134").Value = 2690.4666
$ws3.Range("J134").Value = 13679
$ws3.Range("K134").Value = 8071.399800000001
$ws3.Range("L134").Value = 41037
$ws3.Range("M134").Value = -5536.399800000001
$ws3.Range("N134").Value = -46107

# CRP row 8
$ws4.Range("H8").Value = 500
$ws4.Range("I8").Value = 399
$ws4.Range("J8").Value = 601
$ws4.Range("K8").Value = 399
$ws4.Range("L8").Value = 601
$ws4.Range("N8").Value = -881
$ws4.Range("M8").Value = -259

# CRP row 31
$ws4.Range("H31").Value = 17388.672
$ws4.Range("I31").Value = 1224.5
$ws4.Range("K31").Value = 1224.5
$ws4.Range("M31").Value = -929.5

# CRP row 34
$ws4.Range("H34").Value = 17388.672
$ws4.Range("I34").Value = 1224.5
$ws4.Range("K34").Value = 1224.5
$ws4.Range("M34").Value = -1022.5

# CRP row 41
$ws4.Range("H41").Value = 15744.75
$ws4.Range("I41").Value = 1500
$ws4.Range("K41").Value = 1500
$ws4.Range("M41").Value = -1072

# CRP row 58
$ws4.Range("H58").Value = 3127.3171
$ws4.Range("I58").Value = 3523.25
$ws4.Range("J58").Value = 2274.5386
$ws4.Range("K58").Value = 3523.25
$ws4.Range("L58").Value = 2274.5386
$ws4.Range("M58").Value = -3320.25
$ws4.Range("N58").Value = -2680.5386

# CRP row 122
$ws4.Range("H122").Value = 2284.3333
$ws4.Range("I122").Value = 1898.4667
$ws4.Range("K122").Value = 5695.4001
$ws4.Range("M122").Value = -3245.4001

# CRP row 136
$ws4.Range("H136").Value = 3127.3171
$ws4.Range("I136").Value = 3523.25
$ws4.Range("J136").Value = 2274.5386
$ws4.Range("K136").Value = 10569.75
$ws4.Range("L136").Value = 6823.6158
$ws4.Range("M136").Value = -8019.75
$ws4.Range("N136").Value = -11923.6158

# CUL row 68
$ws5.Range("H68").Value = 872.44446
$ws5.Range("I68").Value = 825
$ws5.Range("J68").Value = 886
$ws5.Range("K68").Value = 2475
$ws5.Range("L68").Value = 2658
$ws5.Range("M68").Value = -1664
$ws5.Range("N68").Value = -4280

# CUL row 71
$ws5.Range("H71").Value = 872.44446
$ws5.Range("I71").Value = 825
$ws5.Range("J71").Value = 886
$ws5.Range("K71").Value = 7425
$ws5.Range("L71").Value = 7974
$ws5.Range("M71").Value = -3369
$ws5.Range("N71").Value = -16086

# CUL row 74
$ws5.Range("H74").Value = 14999.5
$ws5.Range("J74").Value = 14999.5
$ws5.Range("L74").Value = 44998.5
$ws5.Range("N74").Value = -47120.5

# CUL row 77
$ws5.Range("H77").Value = 14999.5
$ws5.Range("J77").Value = 14999.5
$ws5.Range("L77").Value = 134995.5
$ws5.Range("N77").Value = -145603.5

# CUL row 110
$ws5.Range("H110").Value = 14814.571
$ws5.Range("I110").Value = 1881
$ws5.Range("J110").Value = 21999.889
$ws5.Range("K110").Value = 5643
$ws5.Range("L110").Value = 65999.667
$ws5.Range("M110").Value = -1553
$ws5.Range("N110").Value = -74179.667

# CUL row 113
$ws5.Range("H113").Value = 2938.8438
$ws5.Range("J113").Value = 1857.1305
$ws5.Range("L113").Value = 5571.3915
$ws5.Range("N113").Value = -9911.3915

# CUL row 132
$ws5.Range("H132").Value = 2638.353
$ws5.Range("J132").Value = 2952.6365
$ws5.Range("L132").Value = 26573.7285
$ws5.Range("N132").Value = -31633.7285

# CUL row 137
$ws5.Range("H137").Value = 1806
$ws5.Range("I137").Value = 1806
$ws5.Range("J137").Value = 0
$ws5.Range("K137").Value = 5418
$ws5.Range("L137").Value = 0
$ws5.Range("M137").Value = -318
$ws5.Range("N137").ClearContents()

# GSM row 70
$ws6.Range("H70").Value = 6255613.5
$ws6.Range("I70").Value = 8701010
$ws6.Range("K70").Value = 8701010
$ws6.Range("M70").Value = -8700740

# GSM row 73
$ws6.Range("H73").Value = 6255613.5
$ws6.Range("I73").Value = 8701010
$ws6.Range("K73").Value = 8701010
$ws6.Range("M73").Value = -8700074

# GSM row 122
$ws6.Range("H122").Value = 359211.1
$ws6.Range("I122").Value = 557547.25
$ws6.Range("J122").Value = 6613.4443
$ws6.Range("K122").Value = 1672641.75
$ws6.Range("L122").Value = 19840.3329
$ws6.Range("M122").Value = -1670191.75
$ws6.Range("N122").Value = -24740.3329

# GSM row 132
$ws6.Range("H132").Value = 3529.4348
$ws6.Range("I132").Value = 3008.95
$ws6.Range("K132").Value = 9026.849999999999
$ws6.Range("M132").Value = -6496.849999999999

# LTW row 7
$ws7.Range("H7").Value = 4026.182
$ws7.Range("I7").Value = 2706.8667
$ws7.Range("K7").Value = 2706.8667
$ws7.Range("M7").Value = -2594.8667

# LTW row 13
$ws7.Range("H13").Value = 12500
$ws7.Range("I13").Value = 12500
$ws7.Range("J13").Value = 0
$ws7.Range("K13").Value = 12500
$ws7.Range("L13").Value = 0
$ws7.Range("N13").ClearContents()
$ws7.Range("M13").Value = -12360

# LTW row 25
$ws7.Range("H25").Value = 0
$ws7.Range("I25").Value = 0
$ws7.Range("K25").Value = 0
$ws7.Range("M25").ClearContents()

# LTW row 100
$ws7.Range("H100").Value = 3238.5
$ws7.Range("I100").Value = 3154.6191
$ws7.Range("K100").Value = 3154.6191
$ws7.Range("M100").Value = -2613.6191

# LTW row 126
$ws7.Range("H126").Value = 4026.182
$ws7.Range("I126").Value = 2706.8667
$ws7.Range("K126").Value = 8120.6001
$ws7.Range("M126").Value = -5650.6001

# LTW row 132
$ws7.Range("H132").Value = 4469.1304
$ws7.Range("I132").Value = 3670.946
$ws7.Range("K132").Value = 11012.838
$ws7.Range("M132").Value = -8482.838

# WVR row 18
$ws8.Range("H18").Value = 0
$ws8.Range("J18").Value = 0
$ws8.Range("L18").Value = 0
$ws8.Range("N18").ClearContents()

# WVR row 132
$ws8.Range("H132").Value = 34515668
$ws8.Range("I132").Value = 38463560
$ws8.Range("K132").Value = 115390680
$ws8.Range("M132").Value = -115390680
